$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.521.76'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").Value = '3.758.27'
$ws.Range("E3").Value = '  -0.89%  '

$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").Value = "'593.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.80%  '

$ws.Range("D6").Value = "'167.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.85%  '

$ws.Range("D7").Value = '3.755.95'
$ws.Range("E7").Value = '  -0.93%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").Value = "'0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.12%  '

$ws.Range("E10").Value = '  -3.17%  '

$ws.Range("D11").Value = "'6.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.05%  '

$ws.Range("D12").Value = "'0.449"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.17%  '

$ws.Range("D13").Value = "'0.0000260"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -7.40%  '

$ws.Range("D14").Value = "'36.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.22%  '

$ws.Range("D15").Value = '4.392.09'
$ws.Range("E15").Value = '  -0.76%  '

$ws.Range("D16").Value = '3.752.59'
$ws.Range("E16").Value = '  -0.86%  '

$ws.Range("D17").Value = '68.503.43'
$ws.Range("E17").Value = '  +0.52%  '

$ws.Range("E18").Value = '  -3.99%  '

$ws.Range("E19").Value = '  +0.78%  '

$ws.Range("D20").Value = "'7.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.96%  '

$ws.Range("D21").Value = "'10.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.18%  '

$ws.Range("D22").Value = "'465.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.99%  '

$ws.Range("D23").Value = "'0.698"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.26%  '

$ws.Range("D24").Value = "'0.0000148"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.91%  '

$ws.Range("D25").Value = "'84.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.31%  '

$ws.Range("E26").Value = '  -3.30%  '

$ws.Range("D27").Value = "'11.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.85%  '

$ws.Range("D28").Value = "'10.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.04%  '

$ws.Range("E29").Value = '  -0.11%  '

$ws.Range("D30").Value = '3.907.52'
$ws.Range("E30").Value = '  -0.79%  '

$ws.Range("E31").Value = '  -4.83%  '

$ws.Range("D32").Value = "'7.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.83%  '

$ws.Range("D33").Value = "'30.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.11%  '

$ws.Range("D34").Value = "'2.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.83%  '

$ws.Range("D35").Value = "'9.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.17%  '

$ws.Range("D36").Value = "'0.997"
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").Value = '3.715.74'
$ws.Range("E37").Value = '  -0.88%  '

$ws.Range("D38").Value = "'0.100"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.92%  '

$ws.Range("D39").Value = "'3.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -10.04%  '

$ws.Range("E40").Value = '  -1.25%  '

$ws.Range("E41").Value = '  -0.70%  '

$ws.Range("D42").Value = "'5.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.45%  '

$ws.Range("E43").Value = '  +0.17%  '

$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("D45").Value = "'44.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.02%  '

$ws.Range("E46").Value = '  -3.78%  '

$ws.Range("D47").Value = "'46.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.27%  '

$ws.Range("D48").Value = "'1.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.26%  '

$ws.Range("D49").Value = "'8.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.48%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = "'145.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.93%  '

$ws.Range("B51").Value = 'Bittensor'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D51").Value = "'389.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.54%  '
